$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per row (2-9) for columns D, M, N, O, P, R, S
# Row 2
$ws.Range("D2").Value = 44890
$ws.Range("M2").Value = 150
$ws.Range("N2").Value = 13000
$ws.Range("O2").Value = 13000
$ws.Range("P2").Value = 13000
$ws.Range("R2").Value = "La Ligua"
$ws.Range("S2").Value = 2600

# Row 3
$ws.Range("D3").Value = 44890
$ws.Range("M3").Value = 170
$ws.Range("N3").Value = 11000
$ws.Range("O3").Value = 11000
$ws.Range("P3").Value = 11000
$ws.Range("R3").Value = "La Ligua"
$ws.Range("S3").Value = 2200

# Row 4
$ws.Range("D4").Value = 44890
$ws.Range("M4").Value = 150
$ws.Range("N4").Value = 8000
$ws.Range("O4").Value = 8000
$ws.Range("P4").Value = 8000
$ws.Range("R4").Value = "La Ligua"
$ws.Range("S4").Value = 1600

# Row 5
$ws.Range("D5").Value = 44890
$ws.Range("M5").Value = 80
$ws.Range("N5").Value = 7000
$ws.Range("O5").Value = 7000
$ws.Range("P5").Value = 7000
$ws.Range("R5").Value = "La Ligua"
$ws.Range("S5").Value = 1400

# Row 6
$ws.Range("D6").Value = 44908
$ws.Range("M6").Value = 110
$ws.Range("N6").Value = 7000
$ws.Range("O6").Value = 7000
$ws.Range("P6").Value = 7000
$ws.Range("R6").Value = "Provincia de Limarí"
$ws.Range("S6").Value = 1400

# Row 7
$ws.Range("D7").Value = 44908
$ws.Range("M7").Value = 120
$ws.Range("N7").Value = 6000
$ws.Range("O7").Value = 6000
$ws.Range("P7").Value = 6000
$ws.Range("R7").Value = "Provincia de Limarí"
$ws.Range("S7").Value = 1200

# Row 8
$ws.Range("D8").Value = 44908
$ws.Range("M8").Value = 100
$ws.Range("N8").Value = 5000
$ws.Range("O8").Value = 5000
$ws.Range("P8").Value = 5000
$ws.Range("R8").Value = "Provincia de Limarí"
$ws.Range("S8").Value = 1000

# Row 9
$ws.Range("D9").Value = 44908
$ws.Range("M9").Value = 120
$ws.Range("N9").Value = 4000
$ws.Range("O9").Value = 4000
$ws.Range("P9").Value = 4000
$ws.Range("R9").Value = "Provincia de Limarí"
$ws.Range("S9").Value = 800
